# Fix contaminated price/ticker data in columns D:I so every row reflects
# AVAV (AeroVironment) own OHLC prices, shares outstanding, and fixed_ticker,
# instead of data accidentally pulled in from unrelated tickers ("extra files").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
  @{Row=2; D=25.52000045776367; E=27.40999984741211; F=28.5; G=25.43000030517578; H=49932207},
  @{Row=3; D=25.59000015258789; E=25.85000038146973; F=27.06999969482422; G=25.01000022888184; H=49932207},
  @{Row=4; D=26.10000038146973; E=24.04999923706055; F=27; G=22; H=49932207},
  @{Row=5; D=23.03000068664551; E=25.57999992370605; F=26.20999908447266; G=21.86000061035156; H=49932207},
  @{Row=6; D=25.26000022888184; E=24.8799991607666; F=26.42000007629395; G=23.1299991607666; H=49932207},
  @{Row=7; D=28.86000061035156; E=28.80999946594238; F=29.39999961853028; G=27.11000061035156; H=49932207},
  @{Row=8; D=28.44000053405762; E=24.73999977111816; F=30.07999992370605; G=23.85000038146973; H=49932207},
  @{Row=9; D=24.1200008392334; E=28.29999923706055; F=28.79999923706055; G=23.39999961853028; H=49932207},
  @{Row=10; D=26.29999923706055; E=27.02000045776367; F=27.46999931335449; G=25.52000045776367; H=49932207},
  @{Row=11; D=28.6299991607666; E=31.04000091552734; F=31.65999984741211; G=28.1299991607666; H=49932207},
  @{Row=12; D=37.75; E=49.02000045776367; F=49.18000030517578; G=36.70999908447266; H=49932207},
  @{Row=13; D=51.52000045776367; E=45.58000183105469; F=51.72000122070312; G=41.52999877929688; H=49932207},
  @{Row=14; D=51.25; E=49.72000122070312; F=53.11000061035156; G=45.59999847412109; H=49932207},
  @{Row=15; D=54.27000045776367; E=57.86000061035156; F=62.20000076293945; G=49.68999862670898; H=49932207},
  @{Row=16; D=73.81999969482422; E=87.95999908447266; F=88.45999908447266; G=73.5; H=49932207},
  @{Row=17; D=91; E=76.5999984741211; F=103.4599990844727; G=75.2699966430664; H=49932207},
  @{Row=18; D=77.69999694824219; E=79.66999816894531; F=83.5; G=73.72000122070312; H=49932207},
  @{Row=19; D=68.81999969482422; E=64.75; F=70.8499984741211; G=62.5; H=49932207},
  @{Row=20; D=54.84999847412109; E=51.52999877929688; F=55.33000183105469; G=48.61000061035156; H=49932207},
  @{Row=21; D=58.25; E=61.34000015258789; F=63.40000152587891; G=58.02999877929688; H=49932207},
  @{Row=22; D=66.45999908447266; E=51.38999938964844; F=72.69999694824219; G=50.52000045776367; H=49932207},
  @{Row=23; D=58.75; E=70.83000183105469; F=71.06999969482422; G=57.4900016784668; H=49932207},
  @{Row=24; D=77.05999755859375; E=76.38999938964844; F=87; G=75.19999694824219; H=49932207},
  @{Row=25; D=76.9000015258789; E=85.38999938964844; F=89.33999633789062; G=75.12999725341797; H=49932207},
  @{Row=26; D=120.0400009155273; E=110.0800018310547; F=142.2899932861328; G=107.7300033569336; H=49932207},
  @{Row=27; D=110.3499984741211; E=109.629997253418; F=111.7200012207031; G=98.18000030517578; H=49932207},
  @{Row=28; D=101.8499984741211; E=102.3600006103516; F=105.9800033569336; G=96; H=49932207},
  @{Row=29; D=89.75; E=80.76000213623047; F=96.12000274658205; G=79.70999908447266; H=49932207},
  @{Row=30; D=57.34000015258789; E=71.05000305175781; F=71.4800033569336; G=55.34000015258789; H=49932207},
  @{Row=31; D=80.77999877929688; E=91.95999908447266; F=97.1999969482422; G=72.52999877929688; H=49932207},
  @{Row=32; D=85.87000274658203; E=88.6500015258789; F=106.5400009155273; G=85.36000061035156; H=49932207},
  @{Row=33; D=93.33999633789062; E=91.98999786376952; F=94.31999969482422; G=81.44999694824219; H=49932207},
  @{Row=34; D=89.0999984741211; E=85.73999786376953; F=95.5; G=83.55000305175781; H=49932207},
  @{Row=35; D=101.0699996948242; E=93.41000366210938; F=112.3899993896484; G=88.23999786376953; H=49932207},
  @{Row=36; D=95.18000030517578; E=97.02999877929688; F=99.5199966430664; G=91.25; H=49932207},
  @{Row=37; D=115.4700012207031; E=137.6100006103516; F=139.8800048828125; G=114.8600006103516; H=49932207},
  @{Row=38; D=121.9800033569336; E=126.7900009155273; F=128.5; G=119.4700012207031; H=49932207},
  @{Row=39; D=159.1699981689453; E=202.1499938964844; F=205.2799987792969; G=157.1799926757812; H=49932207},
  @{Row=40; D=179.5899963378906; E=203.759994506836; F=211.4400024414062; G=152.8800048828125; H=49932207},
  @{Row=41; D=217.1300048828125; E=194.5; F=236.6000061035156; G=189.259994506836; H=49932207},
  @{Row=42; D=174.8300018310547; E=149.6199951171875; F=188.7899932861328; G=145.9700012207031; H=49932207},
  @{Row=43; D=153.1199951171875; E=178.0299987792969; F=180.7200012207031; G=150.2899932861328; H=49932207},
  @{Row=44; D=259.6099853515625; E=241.3500061035156; F=276.5; G=227.5500030517578; H=49932207}
)

foreach ($item in $rowData) {
  $r = $item.Row
  $ws.Range("D$r").Value = $item.D
  $ws.Range("E$r").Value = $item.E
  $ws.Range("F$r").Value = $item.F
  $ws.Range("G$r").Value = $item.G
  $ws.Range("H$r").Value = $item.H
  $ws.Range("I$r").Value = "AVAV"
}

Write-Host "Updated $($rowData.Count) rows (D:I) to use AVAV-specific data."